$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "DATE" value in D12/D14 was changed to "TIMESTAMP"
$ws.Range("D12").Value = "TIMESTAMP"
$ws.Range("D14").Value = "TIMESTAMP"

# The active selection moved from C7 to D14
$ws.Range("D14").Select()
